$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Relocate the "_GoBack" bookmark.
#    It currently sits in the first paragraph ("Out of 6 5 points").
#    Adding a bookmark with the same name moves/replaces it (bookmark
#    names are unique per document), so simply adding a new "_GoBack"
#    bookmark over the desired span removes the old one automatically.
#    Target span: from the start of the "...best of three..." paragraph
#    through the end of the "...single player..." paragraph.
# ------------------------------------------------------------------
$pBestOfThree = $d.Paragraphs(8)
$pSinglePlayer = $d.Paragraphs(9)
$bmRange = $d.Range($pBestOfThree.Range.Start, $pSinglePlayer.Range.End)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 2. Strike through the obsolete "User stories" items (everything from
#    the "(5 points): ... consistent commits." story through the
#    "(10 points): ... single player ..." story) - they describe the
#    now-removed Round class / its associated stories.
# ------------------------------------------------------------------

# A couple of the paragraphs end up with their sentences broken into
# extra runs (mirrors how Word's own proofing pass re-split them) -
# apply strikethrough to the sub-phrase first so the run boundary is
# created, then restate it over the whole paragraph afterwards so every
# run - old and new - ends up with identical, clean <w:strike/> formatting.
$sub1 = $d.Content
$null = $sub1.Find.Execute("to properly incorporate", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sub1.Font.StrikeThrough = $true

$sub2 = $d.Content
$null = $sub2.Find.Execute("is validated and reobtained if necessary", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sub2.Font.StrikeThrough = $true

For ($i = 3; $i -le 9; $i++) {
    $para = $d.Paragraphs($i)
    $para.Range.Font.StrikeThrough = $true
}

Write-Output "done"
